# Update automatic: dades i banners [2026-02-20 06:36]
# Refresh the DATA_EXTRACCIO timestamp column (E) for each station row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$timestamps = @{
    2 = "2026-02-20 06:34:26"
    3 = "2026-02-20 06:34:29"
    4 = "2026-02-20 06:34:31"
    5 = "2026-02-20 06:34:34"
    6 = "2026-02-20 06:34:36"
    7 = "2026-02-20 06:34:38"
    8 = "2026-02-20 06:34:41"
    9 = "2026-02-20 06:34:43"
    10 = "2026-02-20 06:34:46"
    11 = "2026-02-20 06:34:48"
    12 = "2026-02-20 06:34:51"
    13 = "2026-02-20 06:34:53"
    14 = "2026-02-20 06:34:55"
    15 = "2026-02-20 06:34:58"
    16 = "2026-02-20 06:35:00"
    17 = "2026-02-20 06:35:03"
    18 = "2026-02-20 06:35:05"
    19 = "2026-02-20 06:35:08"
    20 = "2026-02-20 06:35:10"
    21 = "2026-02-20 06:35:13"
    22 = "2026-02-20 06:35:15"
    23 = "2026-02-20 06:35:18"
    24 = "2026-02-20 06:35:20"
    25 = "2026-02-20 06:35:23"
    26 = "2026-02-20 06:35:25"
    27 = "2026-02-20 06:35:27"
    28 = "2026-02-20 06:35:30"
    29 = "2026-02-20 06:35:32"
    30 = "2026-02-20 06:35:35"
    31 = "2026-02-20 06:35:37"
    32 = "2026-02-20 06:35:40"
    33 = "2026-02-20 06:35:42"
    34 = "2026-02-20 06:35:45"
    35 = "2026-02-20 06:35:47"
    36 = "2026-02-20 06:35:50"
    37 = "2026-02-20 06:35:52"
    38 = "2026-02-20 06:35:55"
    39 = "2026-02-20 06:35:57"
    40 = "2026-02-20 06:36:00"
    41 = "2026-02-20 06:36:02"
    42 = "2026-02-20 06:36:05"
    43 = "2026-02-20 06:36:07"
    44 = "2026-02-20 06:36:09"
    45 = "2026-02-20 06:36:12"
    46 = "2026-02-20 06:36:14"
}

foreach ($row in $timestamps.Keys) {
    $ws.Cells.Item($row, 5).Value = $timestamps[$row]
}

